# Add sample data from table 4 to the "ldsurvival-inputs" sheet, and make
# that sheet the active/selected tab (as the commit message + diff show).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ldsurvival-inputs")

# New :t-stage rows (24-30) appended below the existing table.
$ws.Range("A24").Value = ":t-stage"
$ws.Range("B24").Value = ":pT1a"
$ws.Range("D24").Value = 0

$ws.Range("B25").Value = ":pT1b"
$ws.Range("D25").Value = 2

$ws.Range("B26").Value = ":pT2"
$ws.Range("D26").Value = 3

$ws.Range("B27").Value = ":pT3a"
$ws.Range("D27").Value = 4

$ws.Range("B28").Value = ":pT3b"
$ws.Range("D28").Value = 4

$ws.Range("B29").Value = ":pT3c"
$ws.Range("D29").Value = 4

$ws.Range("B30").Value = ":pT4"
$ws.Range("D30").Value = 4

# Make "ldsurvival-inputs" the active sheet/tab, with D30 (the last cell
# entered) as the active selection, matching the saved view state.
$ws.Activate()
$ws.Range("D30").Select()
